# [WIP] wifi check history
#
# The "tbl_checklist_wifi" table schema listing (rows 9-12 of Sheet1)
# gains a new field: "rp_date" (type "varchar"), inserted right after
# "rp_status" and before "remark" -- i.e. a new column is inserted at
# column T, pushing the former T:W columns (remark, user_id,
# approve_status, approver_id) one slot right to U:X.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at T (shifts remark/user_id/approve_status/approver_id right)
$ws.Columns("T:T").Insert()

# New field name/type for the inserted column
$ws.Range("T10").Value = "rp_date"
$ws.Range("T11").Value = "varchar"

# Match the author's final selection/view position
$ws.Range("T12").Select()
